# Generate Report for Handoff
#
# For the six "Ready for handoff" rows (7, 9, 10, 11, 12, 14) that have not yet
# been handed back, this:
#   - stamps a fresh handoff report-generation timestamp
#     (Overview!G and the per-locale "Latest Handoff Datetime" column), and
#   - marks the handoff Priority as "ht".

$wb = $excel.ActiveWorkbook

$rows = 7,9,10,11,12,14

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-19 06:22:02"
}

# --- zh-cn sheet: refresh "Latest Handoff Datetime" (H) and set Priority (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-19 06:21:55"
}

# --- de-de sheet: refresh "Latest Handoff Datetime" (H) and set Priority (E) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-19 06:22:02"
}
